$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.298.50"
$ws.Range("E2").Value = "  +1.55%  "

$ws.Range("D3").Value = "3.381.35"
$ws.Range("E3").Value = "  +0.97%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'579.80"
$ws.Range("E5").Value = "  -0.61%  "

$ws.Range("D6").Value = "'178.23"
$ws.Range("E6").Value = "  +0.55%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("E8").Value = "  +0.50%  "

$ws.Range("E9").Value = "  +7.70%  "

$ws.Range("E10").Value = "  +0.83%  "

$ws.Range("D11").Value = "'48.27"
$ws.Range("E11").Value = "  +0.44%  "

$ws.Range("D12").Value = "'0.0000282"
$ws.Range("E12").Value = "  +3.52%  "

$ws.Range("D13").Value = "'683.68"
$ws.Range("E13").Value = "  -0.80%  "

$ws.Range("E14").Value = "  +1.94%  "

$ws.Range("D15").Value = "3.921.67"
$ws.Range("E15").Value = "  +0.70%  "

$ws.Range("D16").Value = "69.416.04"
$ws.Range("E16").Value = "  +1.49%  "

$ws.Range("E17").Value = "  +0.77%  "

$ws.Range("D18").Value = "3.379.12"
$ws.Range("E18").Value = "  +1.44%  "

$ws.Range("D19").Value = "'17.68"
$ws.Range("E19").Value = "  +1.28%  "

$ws.Range("E20").Value = "  +0.62%  "

$ws.Range("D21").Value = "'0.907"
$ws.Range("E21").Value = "  +1.38%  "

$ws.Range("E22").Value = "  +0.66%  "

$ws.Range("E23").Value = "  -1.88%  "

$ws.Range("D24").Value = "'101.10"
$ws.Range("E24").Value = "  +1.12%  "

$ws.Range("D25").Value = "'3.87"
$ws.Range("E25").Value = "  -0.90%  "

$ws.Range("E26").Value = "  -0.19%  "

$ws.Range("D27").Value = "'9.69"
$ws.Range("E27").Value = "  +1.57%  "

$ws.Range("E28").Value = "  +1.27%  "

$ws.Range("D29").Value = "'8.70"
$ws.Range("E29").Value = "  +2.53%  "

$ws.Range("D30").Value = "'6.91"
$ws.Range("E30").Value = "  -0.23%  "

$ws.Range("E32").Value = "  -0.53%  "

$ws.Range("D33").Value = "'548.81"
$ws.Range("E33").Value = "  -2.12%  "

$ws.Range("E34").Value = "  -0.17%  "

$ws.Range("E35").Value = "  +0.27%  "

$ws.Range("E36").Value = "  -0.08%  "

$ws.Range("D37").Value = "3.600.14"
$ws.Range("E37").Value = "  -2.74%  "

$ws.Range("E38").Value = "  +3.17%  "

$ws.Range("D39").Value = "'35.26"
$ws.Range("E39").Value = "  +1.72%  "

$ws.Range("D40").Value = "0.0₃0737"
$ws.Range("E40").Value = "  +9.70%  "

$ws.Range("E41").Value = "  +4.65%  "

$ws.Range("E42").Value = "  +3.46%  "

$ws.Range("D43").Value = "'3.38"
$ws.Range("E43").Value = "  +3.49%  "

$ws.Range("E44").Value = "  +2.97%  "

$ws.Range("E45").Value = "  -0.29%  "

$ws.Range("E46").Value = "  +0.36%  "

$ws.Range("E47").Value = "  -0.06%  "

$ws.Range("E48").Value = "  +3.57%  "

$ws.Range("E49").Value = "  -0.24%  "

$ws.Range("D50").Value = "'129.51"
$ws.Range("E50").Value = "  -0.98%  "

$ws.Range("D51").Value = "'2.58"
$ws.Range("E51").Value = "  +0.72%  "
